$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.493.25'
$ws.Range('E2').Value = '  -1.97%  '
$ws.Range('D3').Value = '3.689.35'
$ws.Range('E3').Value = '  -2.99%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '681.69'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.32%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '162.81'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.25%  '
$ws.Range('D7').Value = '3.688.57'
$ws.Range('E7').Value = '  -2.97%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -4.22%  '
$ws.Range('E10').Value = '  -7.58%  '
$ws.Range('E11').Value = '  -3.09%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.447'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.48%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000239'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -4.56%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '33.65'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -5.94%  '
$ws.Range('D15').Value = '4.313.04'
$ws.Range('D16').Value = '3.686.31'
$ws.Range('E16').Value = '  -2.31%  '
$ws.Range('D17').Value = '69.538.73'
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.112'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.29%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '16.29'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -6.18%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.63'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -6.48%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '483.23'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.09%  '
$ws.Range('E22').Value = '  -6.84%  '
$ws.Range('E23').Value = '  -7.53%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '80.35'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.32%  '
$ws.Range('D25').Value = '3.834.51'
$ws.Range('E25').Value = '  -3.00%  '
$ws.Range('E26').Value = '  -8.27%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.51'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.20%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.56'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -7.02%  '
$ws.Range('E30').Value = '  -8.08%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.72'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -10.28%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.90'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -5.49%  '
$ws.Range('E33').Value = '  -7.28%  '
$ws.Range('E34').Value = '  -2.43%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '27.16'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -6.22%  '
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('D37').Value = '3.659.49'
$ws.Range('E37').Value = '  -2.89%  '
$ws.Range('E38').Value = '  -5.54%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.35'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +6.89%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0941'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -6.66%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.24'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.24%  '
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('E44').Value = '  -6.11%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '161.41'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.20%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '48.35'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.40%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.86'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -11.40%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '30.20'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +7.53%  '
$ws.Range('E49').Value = '  -7.59%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.35'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.52%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '392.54'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -6.18%  '
